$d = $word.ActiveDocument
$q = [char]34

# --- 1. Update the "generated on" sentence (date, time, runner id). Do this
#     before touching the title date, since this sentence also contains the
#     old date string "2021-08-25" too and we want to target it precisely
#     via a full-sentence match first. Setting .Text directly (rather than
#     relying on Find.Execute's Replace parameter) avoids Word's
#     smart-quote autocorrect mangling the straight quotes around
#     "develop".
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "This document was generated on 2021-08-25, 09:20:42 with the Automatic Report Generator (ARG) version " + $q + "develop" + $q + " on the Linux system runner-fa6cab46-project-18732201-concurrent-0.",
    $true)
if ($found1) {
    $r1.Text = "This document was generated on 2021-08-26, 08:26:23 with the Automatic Report Generator (ARG) version " + $q + "develop" + $q + " on the Linux system runner-ed2dce3a-project-18732201-concurrent-0."
}

# --- 2. Update the title date (the remaining "2021-08-25" occurrence, in
#     the title block).
$r2 = $d.Content
$found2 = $r2.Find.Execute("2021-08-25", $true)
if ($found2) {
    $r2.Text = "2021-08-26"
}

# --- 3. Add a new paragraph after "Just another string: ARG's documentation
#     (Chapter 2)" containing a hyperlink (re-using the same target as the
#     existing ARG's-documentation link, rId11) whose display text is
#     "ARG's documentation but no string before", followed by a plain
#     " (Chapter 2)" run - mirroring the paragraph above it, but without the
#     leading "Just another string: " text run. We build the exact OOXML
#     ourselves (instead of using Hyperlinks.Add / typing + Font formatting)
#     so the run keeps direct character formatting
#     (w:color + w:u, matching the sibling hyperlink already in the
#     document) instead of an auto-generated "Hyperlink" character style.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*ARG's documentation (Chapter 2)*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:hyperlink r:id="rId11"><w:r><w:rPr><w:color w:val="0000FF"/><w:u w:val="single"/></w:rPr><w:t>ARG's documentation but no string before</w:t></w:r></w:hyperlink><w:r><w:rPr/><w:t xml:space="preserve"> (Chapter 2)</w:t></w:r></w:p>
"@
    [void]$newPara.Range.InsertXML($xml)
}
